$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 49

$ws.Cells.Item($row, 1).Value = "XDR4AZ"
$ws.Cells.Item($row, 2).Value = "Kit de engranajes para unidad de revelado RICOH"
$ws.Cells.Item($row, 3).Value = "Aficio 1515, MP161 MP171 MP175 MP201 MP301"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 80000
$ws.Cells.Item($row, 6).Value = 9
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E49-D49)*G49"
$ws.Cells.Item($row, 9).Formula = "=D49*F49"
$ws.Cells.Item($row, 10).Value = 0
